$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "SELECT DISTINCT`n     std.dbgap_accession AS `"dbGaP Accession`",`n     std.study_name AS `"Study Name`"`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`"`nWHERE `n    std.dbgap_accession = 'phs002518' AND dgn.diagnosis = '8000/0 : Neoplasm, benign'"
$ws.Range("C2").Value = "SELECT `n    COUNT(DISTINCT dgn.diagnosis) AS Diagnoses,`n    COUNT(DISTINCT prt.participant_id) AS Participants,`n    COUNT(DISTINCT std.study_id) AS Studies`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`"`nWHERE `n   std.dbgap_accession = 'phs002518' AND dgn.diagnosis = '8000/0 : Neoplasm, benign';"
$ws.Range("B3").Value = "SELECT DISTINCT`n    prt.participant_id AS `"Participant Id`",`n    prt.race AS `"Race`",`n    prt.sex_at_birth AS `"Sex at Birth`",`n    std.dbgap_accession AS `"dbGaP Accession`"`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`"`nWHERE `n    std.dbgap_accession = 'phs002518' AND dgn.diagnosis = '8000/0 : Neoplasm, benign'`nORDER BY `n    prt.participant_id ASC`nLIMIT 100;"
$ws.Range("B4").Value = "SELECT DISTINCT`n    prt.participant_id AS `"Participant Id`",`n    dgn.diagnosis_id AS `"Diagnosis Id`",`n    dgn.diagnosis AS `"Diagnosis`",`n    dgn.diagnosis_classification_system AS `"Diagnosis Classification System`",`n    dgn.diagnosis_basis AS `"Diagnosis Basis`",`n    dgn.tumor_classification AS `"Tumor Classification`",`n    dgn.anatomic_site AS `"Anatomic Site`",`n    CASE `n    WHEN dgn.age_at_diagnosis = -999 THEN 'Not Reported'`n    WHEN dgn.age_at_diagnosis >= 1000 THEN `n        substr(dgn.age_at_diagnosis, 1, length(dgn.age_at_diagnosis) - 3) || ',' || substr(dgn.age_at_diagnosis, -3)`n    ELSE `n        dgn.age_at_diagnosis `nEND AS `"Age at Diagnosis (days)`",`n    std.dbgap_accession AS `"dbGaP Accession`"`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`"`nWHERE `n    std.dbgap_accession = 'phs002518' AND dgn.diagnosis = '8000/0 : Neoplasm, benign' AND dgn.diagnosis_id IS NOT NULL`nORDER BY `n    prt.participant_id ASC`nLIMIT 100;"
$ws.Range("B6").Value = "SELECT DISTINCT`n    prt.participant_id AS `"Participant Id`",`n    trr.treatment_response_id AS `"Treatment Response Id`",`n    trr.response AS `"Response`",`n    CASE `n        WHEN trr.age_at_response = -999 THEN 'Not Reported'`n        WHEN trr.age_at_response >= 1000 THEN `n            substr(trr.age_at_response, 1, length(trr.age_at_response) - 3) || ',' || substr(trr.age_at_response, -3)`n        ELSE `n            trr.age_at_response `n    END AS `"Age at Response`",`n    trr.response_category AS `"Response Category`",`n    trr.response_system AS `"Response System`",`n    std.dbgap_accession AS `"dbGaP Accession`"`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`"`nWHERE `n    std.dbgap_accession = 'phs002518' AND dgn.diagnosis = '8000/0 : Neoplasm, benign'`nORDER BY `n    prt.participant_id ASC`nLIMIT 100;"
$ws.Range("B7").Value = "SELECT DISTINCT`n    prt.participant_id AS `"Participant Id`",`n    srv.survival_id AS `"Survival Id`",`n    srv.last_known_survival_status AS `"Last Known Survival Status`",`n    CASE `n    WHEN srv.age_at_last_known_survival_status = -999 THEN 'Not Reported'`n    WHEN srv.age_at_last_known_survival_status >= 1000 THEN `n        substr(srv.age_at_last_known_survival_status, 1, length(srv.age_at_last_known_survival_status) - 3) || ',' || substr(srv.age_at_last_known_survival_status, -3)`n    ELSE `n        srv.age_at_last_known_survival_status `nEND AS `"Age at Last Known Survival Status`",`n    srv.first_event AS `"First Event`",`n    srv.cause_of_death AS `"Cause of Death`",`n    std.dbgap_accession AS `"dbGaP Accession`"`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`"`nWHERE `n    std.dbgap_accession = 'phs002518' AND dgn.diagnosis = '8000/0 : Neoplasm, benign' AND srv.survival_id IS NOT NULL`nORDER BY `n    prt.participant_id ASC`nLIMIT 100;"
